# More Customization & Better Output
#
# Updates the evaluation metrics for the "Principal Component Coordinates"
# method (row 2) to reflect a re-run of the color-matching benchmark.
# The values are stored as text (shared strings), so a helper cell
# formatted as Text is used to push numeric-looking strings into the
# target cells without Excel reinterpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Row, $Col, $Text)

    $helper = $Worksheet.Cells.Item(1048576, 16384)
    $helper.NumberFormat = "@"
    $helper.Value = $Text
    $helper.Copy()
    $Worksheet.Cells.Item($Row, $Col).PasteSpecial(-4163)
    $helper.Clear()
}

Set-TextValue $ws 2 4 "130"
Set-TextValue $ws 2 5 "0.15662"
Set-TextValue $ws 2 6 "9.52117"
Set-TextValue $ws 2 7 "0.13538"
Set-TextValue $ws 2 8 "0.98369"
Set-TextValue $ws 2 13 "0.00217"

$excel.CutCopyMode = $false
